$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 6; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $d -replace "LibSrc-Genomi_", "LibSrc-Genomic_"
    $ws.Cells.Item($r, 5).Value2 = $e -replace "LibSrc-Genomi_", "LibSrc-Genomic_"
}

$ws.Range("D16").Select() | Out-Null
